$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns AD, AE, AF -> Wins, Losses, Ties
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# copy header style from an existing header cell (e.g. AC1) so formatting matches
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Fill team record data for every data row (2 through 45)
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 30).Value = 73   # AD
    $ws.Cells.Item($r, 31).Value = 89   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
